# Updates the cryptos worksheet with refreshed Price (D) and Volume(1h) (E)
# figures, and fixes the ordering of two coin-pairs whose rows were swapped
# (Fetch.AI / InternetComputer(DFINITY) at rows 25-26, and Mantle / PEPE at
# rows 35-36). Cells that hold values which look like plain numbers are
# forced back to Text format first so Excel keeps them as literal strings
# (e.g. "1.00" instead of being auto-converted to the number 1), matching
# the original inline-string formatting of the Price column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.197.53"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").Value = "3.003.19"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.35"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.39"
$ws.Range("E6").Value = "  -6.12%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").Value = "2.998.91"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.42"
$ws.Range("E14").Value = "  -6.93%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "3.499.71"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.12"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "62.196.52"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").Value = "3.007.55"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.12"
$ws.Range("E20").Value = "  -5.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.43"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.61"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.30"
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -10.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  -5.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.61"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.99"
$ws.Range("E31").Value = "  -6.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.43"
$ws.Range("E32").Value = "  +3.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.07"
$ws.Range("E33").Value = "  -7.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0795"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").Value = "  -4.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.76"
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("E38").Value = "  -5.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.22"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.25"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  -12.13%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "391.78"
$ws.Range("E43").Value = "  -10.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0357"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.270"
$ws.Range("E45").Value = "  -7.12%  "
$ws.Range("D46").Value = "2.722.39"
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.22"
$ws.Range("E47").Value = "  -7.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.59"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.21"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -1.05%  "

